$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# EDUCATION section: the "Bachelors of Arts in Spanish and English"
# entry (heading, its "{UNIVERSITY}" line, and the blank spacer
# paragraph right after it) is removed entirely, leaving "Associates of
# Arts" as the sole remaining education entry.
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("Bachelors of Arts in Spanish and English", $true, $false,
                     $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$idx1 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pr = $d.Paragraphs.Item($i).Range
    if ($pr.Start -le $find1.Start -and $pr.End -gt $find1.Start) {
        $idx1 = $i
        break
    }
}

$start1 = $d.Paragraphs.Item($idx1).Range.Start
$end1 = $d.Paragraphs.Item($idx1 + 2).Range.End
$d.Range($start1, $end1).Delete()

# ---------------------------------------------------------------------
# SKILLS section: drop the "Proficient in Microsoft Word." and
# "Bilingual in Spanish and English." bullets entirely, then re-plant
# the "_GoBack" bookmark at the start of the next bullet ("Ability to
# adapt to new environments and learn quickly.").
# ---------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("Proficient in Microsoft Word.", $true, $false,
                     $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$idx2 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pr = $d.Paragraphs.Item($i).Range
    if ($pr.Start -le $find2.Start -and $pr.End -gt $find2.Start) {
        $idx2 = $i
        break
    }
}

$start2 = $d.Paragraphs.Item($idx2).Range.Start
$end2 = $d.Paragraphs.Item($idx2 + 1).Range.End
$d.Range($start2, $end2).Delete()

$nextPara = $d.Paragraphs.Item($idx2)
$bmStart = $nextPara.Range.Start
$d.Bookmarks.Add("_GoBack", $d.Range($bmStart, $bmStart))
